$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("H6").Value = 3.15
$ws.Range("I6").Value = 3.1
$ws.Range("L6").Value = 3.5
$ws.Range("W6").Value = 7.5
$ws.Range("X6").Value = 10.75
$ws.Range("Z6").Value = 22
$ws.Range("AA6").Value = 18.5
$ws.Range("AB6").Value = 28
$ws.Range("AG6").Value = 9.75
$ws.Range("AH6").Value = 17
$ws.Range("AK6").Value = 26
$ws.Range("AL6").Value = 32
$ws.Range("AM6").Value = 450
$ws.Range("AP6").Value = 18.5
$ws.Range("AR6").Value = 75
$ws.Range("AT6").Value = 2.62
$ws.Range("AU6").Value = 6.6
$ws.Range("AW6").Value = 5.1
$ws.Range("AX6").Value = 16
$ws.Range("AY6").Value = 21
$ws.Range("AZ6").Value = 75
$ws.Range("BA6").Value = 100

# Row 7
$ws.Range("S7").Value = 1.28
$ws.Range("T7").Value = 3.46

# Row 10
$ws.Range("G10").Value = 3.9
$ws.Range("H10").Value = 3.3
$ws.Range("I10").Value = 1.85
$ws.Range("J10").Value = 4.75
$ws.Range("K10").Value = 2.05
$ws.Range("L10").Value = 2.6
$ws.Range("M10").Value = 1.07
$ws.Range("N10").Value = 8.5
$ws.Range("U10").Value = 2
$ws.Range("V10").Value = 1.73
$ws.Range("W10").Value = 9.5
$ws.Range("X10").Value = 19
$ws.Range("Z10").Value = 41
$ws.Range("AA10").Value = 34
$ws.Range("AB10").Value = 41
$ws.Range("AC10").Value = 8.5
$ws.Range("AD10").Value = 6.5
$ws.Range("AH10").Value = 8
$ws.Range("AJ10").Value = 15
$ws.Range("AN10").Value = 6
$ws.Range("AO10").Value = 23
$ws.Range("AP10").Value = 34
$ws.Range("AQ10").Value = 81
$ws.Range("AS10").Value = 301
$ws.Range("AU10").Value = 8.5
$ws.Range("AW10").Value = 3.75
$ws.Range("AX10").Value = 11
$ws.Range("AZ10").Value = 41

# Row 13
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 3.1
$ws.Range("I13").Value = 3.6
$ws.Range("J13").Value = 2.75
$ws.Range("W13").Value = 7
$ws.Range("AD13").Value = 6

# Row 18
$ws.Range("H18").Value = 3
$ws.Range("M18").Value = 1.11
$ws.Range("N18").Value = 6.5
$ws.Range("S18").Value = 1.57
$ws.Range("T18").Value = 2.25
$ws.Range("U18").Value = 2.2
$ws.Range("V18").Value = 1.62
$ws.Range("W18").Value = 5.5
$ws.Range("Z18").Value = 19
$ws.Range("AF18").Value = 81
$ws.Range("AG18").Value = 8.5
$ws.Range("AN18").Value = 4
$ws.Range("AP18").Value = 29
$ws.Range("AR18").Value = 81
$ws.Range("AT18").Value = 2.25
$ws.Range("AU18").Value = 9.5
$ws.Range("AV18").Value = 81
$ws.Range("AY18").Value = 41

# Row 23
$ws.Range("L23").Value = 2.88

Write-Host "Applied odds updates for rows 6, 7, 10, 13, 18, 23"
